# Sara-Alert-Format-Invalid-Monitorees.xlsx
#
# Commit: "add new race options to import/export, reorder race options for
# consistency, update excel files"
#
# This fixture's header row (row 1, columns A:CX) already enumerates every
# import/export column; the only semantic change that lands in THIS file is
# the "reorder race options for consistency" part: the two trailing race
# header cells - "Race Unknown" (CV1) and "Race Other" (CW1) - swap places.
#
# The workbook's dimension stays A1:CX1 (no columns are inserted here), so we
# just swap the two header values in place, preserving each cell's existing
# style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$raceUnknown = $ws.Range("CV1").Value2
$raceOther   = $ws.Range("CW1").Value2

$ws.Range("CV1").Value = $raceOther
$ws.Range("CW1").Value = $raceUnknown

# Best-effort reproduction of the saved selection/viewport: the author's
# last selection before saving was the range BA9:CJ9.
$ws.Range("BA9:CJ9").Select()
